# Append the 2025-02-16 12:22:53 resale-number row to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 72

# Columns A-D hold plain text (date/time/weekday/week strings), not numbers
# or Excel-native dates/times. Temporarily force text formatting so values
# such as "2025-02-16" or "07" are not reinterpreted as a date serial or a
# number, then restore the default "Normal" style so no explicit cell style
# is left behind (matching the rest of the data rows).
$textRange = $ws.Range("A$row`:D$row")
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-02-16"
$ws.Cells.Item($row, 2).Value = "12:22:53"
$ws.Cells.Item($row, 3).Value = "Sunday"
$ws.Cells.Item($row, 4).Value = "07"

$textRange.Style = "Normal"

$ws.Cells.Item($row, 5).Value = 128274
$ws.Cells.Item($row, 6).Value = 140216
$ws.Cells.Item($row, 7).Value = 170642
$ws.Cells.Item($row, 8).Value = 159417
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 145226
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192307
$ws.Cells.Item($row, 14).Value = 115226
$ws.Cells.Item($row, 15).Value = 45165
$ws.Cells.Item($row, 16).Value = 28781
$ws.Cells.Item($row, 17).Value = 66143
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 46321
$ws.Cells.Item($row, 20).Value = -1
